# Apply 2025-11-17 data updates to violent-crime-full-year workbook
# Updates the 2025 (column L) figures - and, where the 2018 figure (column E)
# was also revised, those cells too - across the Citywide Totals, By Neighborhood,
# and individual neighborhood sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item('Citywide Totals')
$ws.Range('L2').Value = 5878
$ws.Range('L3').Value = 6401
$ws.Range('E4').Value = 2064
$ws.Range('L4').Value = 1576
$ws.Range('L5').Value = 383
$ws.Range('L6').Value = 5262
$ws.Range('E7').Value = 26069
$ws.Range('L7').Value = 19500

$ws = $wb.Worksheets.Item('Logan Square')
$ws.Range('L2').Value = 64
$ws.Range('L6').Value = 72
$ws.Range('L7').Value = 214

$ws = $wb.Worksheets.Item('Austin')
$ws.Range('L2').Value = 382
$ws.Range('L3').Value = 456
$ws.Range('L7').Value = 1284

$ws = $wb.Worksheets.Item('Garfield Park')
$ws.Range('L6').Value = 254
$ws.Range('L7').Value = 882

$ws = $wb.Worksheets.Item('West Pullman')
$ws.Range('L6').Value = 63
$ws.Range('L7').Value = 274

$ws = $wb.Worksheets.Item('Grand Crossing')
$ws.Range('L2').Value = 223
$ws.Range('L3').Value = 263
$ws.Range('L7').Value = 744

$ws = $wb.Worksheets.Item('New City')
$ws.Range('L2').Value = 137
$ws.Range('L3').Value = 124
$ws.Range('L6').Value = 92
$ws.Range('L7').Value = 379

$ws = $wb.Worksheets.Item('Woodlawn')
$ws.Range('L3').Value = 140
$ws.Range('L7').Value = 339

$ws = $wb.Worksheets.Item('Fuller Park')
$ws.Range('L2').Value = 24
$ws.Range('L7').Value = 87

$ws = $wb.Worksheets.Item('By Neighborhood')
$ws.Range('L2').Value = 166
$ws.Range('L7').Value = 628
$ws.Range('L8').Value = 1284
$ws.Range('L11').Value = 323
$ws.Range('L12').Value = 45
$ws.Range('L15').Value = 157
$ws.Range('L19').Value = 533
$ws.Range('L20').Value = 493
$ws.Range('L22').Value = 61
$ws.Range('L24').Value = 56
$ws.Range('L30').Value = 87
$ws.Range('L31').Value = 195
$ws.Range('L33').Value = 882
$ws.Range('L34').Value = 110
$ws.Range('L37').Value = 744
$ws.Range('L45').Value = 36
$ws.Range('L48').Value = 259
$ws.Range('L51').Value = 247
$ws.Range('L52').Value = 405
$ws.Range('L53').Value = 214
$ws.Range('L54').Value = 424
$ws.Range('L55').Value = 201
$ws.Range('L57').Value = 67
$ws.Range('L61').Value = 20
$ws.Range('E63').Value = 395
$ws.Range('L63').Value = 56
$ws.Range('L65').Value = 379
$ws.Range('L67').Value = 672
$ws.Range('L76').Value = 297
$ws.Range('L78').Value = 248
$ws.Range('L79').Value = 543
$ws.Range('L85').Value = 964
$ws.Range('L86').Value = 129
$ws.Range('L88').Value = 207
$ws.Range('L90').Value = 204
$ws.Range('L91').Value = 262
$ws.Range('L92').Value = 59
$ws.Range('L94').Value = 239
$ws.Range('L95').Value = 274
$ws.Range('L98').Value = 103
$ws.Range('L99').Value = 339
$ws.Range('L100').Value = 35
$ws.Range('E101').Value = 26069
$ws.Range('L101').Value = 19500

$ws = $wb.Worksheets.Item('Gage Park')
$ws.Range('L2').Value = 79
$ws.Range('L3').Value = 51
$ws.Range('L7').Value = 195

$ws = $wb.Worksheets.Item('North Lawndale')
$ws.Range('L3').Value = 261
$ws.Range('L6').Value = 155
$ws.Range('L7').Value = 672

$ws = $wb.Worksheets.Item('Loop')
$ws.Range('L4').Value = 35
$ws.Range('L6').Value = 208
$ws.Range('L7').Value = 424

$ws = $wb.Worksheets.Item('Lake View')
$ws.Range('L2').Value = 39
$ws.Range('L6').Value = 102
$ws.Range('L7').Value = 259

$ws = $wb.Worksheets.Item('Chatham')
$ws.Range('L2').Value = 188
$ws.Range('L3').Value = 167
$ws.Range('L7').Value = 533

$ws = $wb.Worksheets.Item('River North')
$ws.Range('L2').Value = 62
$ws.Range('L6').Value = 134
$ws.Range('L7').Value = 297

$ws = $wb.Worksheets.Item('Rogers Park')
$ws.Range('L3').Value = 84
$ws.Range('L7').Value = 248

$ws = $wb.Worksheets.Item('Lower West Side')
$ws.Range('L2').Value = 61
$ws.Range('L7').Value = 201

$ws = $wb.Worksheets.Item('Dunning')
$ws.Range('L6').Value = 13
$ws.Range('L7').Value = 56

$ws = $wb.Worksheets.Item('Washington Park')
$ws.Range('L3').Value = 122
$ws.Range('L7').Value = 262

$ws = $wb.Worksheets.Item('Roseland')
$ws.Range('L2').Value = 172
$ws.Range('L3').Value = 174
$ws.Range('L6').Value = 145
$ws.Range('L7').Value = 543

$ws = $wb.Worksheets.Item('Chicago Lawn')
$ws.Range('L2').Value = 154
$ws.Range('L3').Value = 170
$ws.Range('L4').Value = 47
$ws.Range('L7').Value = 493

$ws = $wb.Worksheets.Item('Wrigleyville')
$ws.Range('L6').Value = 21
$ws.Range('L7').Value = 35

$ws = $wb.Worksheets.Item('Auburn Gresham')
$ws.Range('L2').Value = 213
$ws.Range('L7').Value = 628

$ws = $wb.Worksheets.Item('Garfield Ridge')
$ws.Range('L2').Value = 38
$ws.Range('L7').Value = 110

$ws = $wb.Worksheets.Item('West Loop')
$ws.Range('L5').Value = 6
$ws.Range('L7').Value = 239

$ws = $wb.Worksheets.Item('Brighton Park')
$ws.Range('L2').Value = 59
$ws.Range('L7').Value = 157

$ws = $wb.Worksheets.Item('Wicker Park')
$ws.Range('L6').Value = 48
$ws.Range('L7').Value = 103

$ws = $wb.Worksheets.Item('Belmont Cragin')
$ws.Range('L2').Value = 120
$ws.Range('L7').Value = 323

$ws = $wb.Worksheets.Item('Albany Park')
$ws.Range('L2').Value = 55
$ws.Range('L3').Value = 53
$ws.Range('L7').Value = 166

$ws = $wb.Worksheets.Item('West Elsdon')
$ws.Range('L2').Value = 25
$ws.Range('L7').Value = 59

$ws = $wb.Worksheets.Item('United Center')
$ws.Range('L2').Value = 63
$ws.Range('L3').Value = 70
$ws.Range('L7').Value = 207

$ws = $wb.Worksheets.Item('Streeterville')
$ws.Range('L6').Value = 16
$ws.Range('L7').Value = 129

$ws = $wb.Worksheets.Item('Washington Heights')
$ws.Range('L2').Value = 67
$ws.Range('L3').Value = 57
$ws.Range('L7').Value = 204

$ws = $wb.Worksheets.Item('Little Italy, UIC')
$ws.Range('L4').Value = 36
$ws.Range('L7').Value = 247

$ws = $wb.Worksheets.Item('Mckinley Park')
$ws.Range('L2').Value = 21
$ws.Range('L7').Value = 67

$ws = $wb.Worksheets.Item('South Shore')
$ws.Range('L3').Value = 396
$ws.Range('L6').Value = 201
$ws.Range('L7').Value = 964

$ws = $wb.Worksheets.Item('Clearing')
$ws.Range('L3').Value = 22
$ws.Range('L7').Value = 61

$ws = $wb.Worksheets.Item('Jackson Park')
$ws.Range('L6').Value = 9
$ws.Range('L7').Value = 36

$ws = $wb.Worksheets.Item('Little Village')
$ws.Range('L6').Value = 109
$ws.Range('L7').Value = 405

$ws = $wb.Worksheets.Item('Beverly')
$ws.Range('L2').Value = 13
$ws.Range('L7').Value = 45

$ws = $wb.Worksheets.Item('Mount Greenwood')
$ws.Range('L2').Value = 6
$ws.Range('L7').Value = 20
